$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 - species data swapped from old row 28
$ws.Range("A26").Value = 112093190
$ws.Range("B26").Value = 85210
$ws.Range("E26").Value = 3624
$ws.Range("F26").Value = "Strimspindling"
$ws.Range("G26").Value = "Cortinarius glaucopus"
$ws.Range("H26").Value = "(Schaeff. : Fr.) Fr."

# Row 27 - species data swapped from old row 29, plus location + comment updates
$ws.Range("A27").Value = 112093171
$ws.Range("B27").Value = 88909
$ws.Range("D27").Value = "VU"
$ws.Range("E27").Value = 720
$ws.Range("F27").Value = "Violgubbe"
$ws.Range("G27").Value = "Gomphus clavatus"
$ws.Range("H27").Value = "(Pers.) Gray"
$ws.Range("P27").Value = "Bladsätra, Upl"
$ws.Range("Q27").Value = 639204.9761395331
$ws.Range("R27").Value = 6701015.582563667
$ws.Range("AC27").Value = "Barkborredödat bestånd som även kantar mot ett stort kalhygge"
$ws.Range("AD27").Value = $true

# Row 28 - species data swapped from old row 26
$ws.Range("A28").Value = 112093192
$ws.Range("B28").Value = 90687
$ws.Range("E28").Value = 5964
$ws.Range("F28").Value = "Fjällig taggsvamp s.str."
$ws.Range("G28").Value = "Sarcodon imbricatus s.str."
$ws.Range("H28").Value = "(L.:Fr.) P.Karst."

# Row 29 - species data swapped from old row 27, plus location + comment updates
$ws.Range("A29").Value = 112093193
$ws.Range("B29").Value = 89183
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 3215
$ws.Range("F29").Value = "Rödgul trumpetsvamp"
$ws.Range("G29").Value = "Craterellus lutescens"
$ws.Range("H29").Value = "(Fr.) Fr."
$ws.Range("Q29").Value = 639179.9128251362
$ws.Range("R29").Value = 6701165.391882338
$ws.Range("AC29").ClearContents()
$ws.Range("AD29").Value = $false
